$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows ---
# Row 64
$ws.Cells.Item(64, 11).Value = "Finished"
$ws.Cells.Item(64, 20).Value = "2024-07-12 16:48"

# Row 68
$ws.Cells.Item(68, 11).Value = "Paused"
$ws.Cells.Item(68, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(68, 15).Value = "Fin del día automático"

# Row 69
$ws.Cells.Item(69, 11).Value = "Finished"
$ws.Cells.Item(69, 14).Value = "2024-07-12 07:51"
$ws.Cells.Item(69, 15).Value = "Fin del día"
$ws.Cells.Item(69, 16).Value = "2024-07-12 07:51"
$ws.Cells.Item(69, 20).Value = "2024-07-12 09:02"

# Row 71
$ws.Cells.Item(71, 11).Value = "Paused"
$ws.Cells.Item(71, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(71, 15).Value = "Fin del día automático"

# Row 77
$ws.Cells.Item(77, 11).Value = "Paused"
$ws.Cells.Item(77, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(77, 15).Value = "Fin del día automático"

# Row 79
$ws.Cells.Item(79, 11).Value = "Finished"
$ws.Cells.Item(79, 16).Value = "2024-07-12 12:51"
$ws.Cells.Item(79, 17).Value = "2024-07-12 12:52"
$ws.Cells.Item(79, 18).Value = "Falta de materiales: Corrugado"
$ws.Cells.Item(79, 20).Value = "2024-07-12 12:52"

# Row 80
$ws.Cells.Item(80, 16).Value = "2024-07-13 15:56"
$ws.Cells.Item(80, 17).Value = "2024-07-13 15:57"
$ws.Cells.Item(80, 18).Value = "Fin del día automático"

# Row 81
$ws.Cells.Item(81, 11).Value = "Paused"
$ws.Cells.Item(81, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(81, 15).Value = "Fin del día automático"

# Row 82
$ws.Cells.Item(82, 11).Value = "Paused"
$ws.Cells.Item(82, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(82, 15).Value = "Fin del día automático"

# Row 83
$ws.Cells.Item(83, 11).Value = "Paused"
$ws.Cells.Item(83, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(83, 15).Value = "Fin del día automático"

# Row 84
$ws.Cells.Item(84, 11).Value = "Paused"
$ws.Cells.Item(84, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(84, 15).Value = "Fin del día automático"

# Row 86
$ws.Cells.Item(86, 11).Value = "Paused"
$ws.Cells.Item(86, 14).Value = "2024-07-12 07:52"
$ws.Cells.Item(86, 15).Value = "Pausa para almorzar"

# Row 88
$ws.Cells.Item(88, 11).Value = "Finished"
$ws.Cells.Item(88, 20).Value = "2024-07-12 07:50"

# Row 91
$ws.Cells.Item(91, 11).Value = "Paused"
$ws.Cells.Item(91, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(91, 15).Value = "Fin del día automático"

# Row 92
$ws.Cells.Item(92, 11).Value = "Paused"
$ws.Cells.Item(92, 14).Value = "2024-07-12 12:51"
$ws.Cells.Item(92, 15).Value = "Pausa para almorzar"
$ws.Cells.Item(92, 16).Value = ""
$ws.Cells.Item(92, 17).Value = ""
$ws.Cells.Item(92, 18).Value = ""
$ws.Cells.Item(92, 19).Value = ""
$ws.Cells.Item(92, 20).Value = ""

# --- New rows 93-101 ---
# Row 93
$ws.Cells.Item(93, 1).Value = "8218092e-4e36-4b61-960c-0e971b07466b"
$ws.Cells.Item(93, 2).Value = "2024-07-12 07:50"
$ws.Cells.Item(93, 3).Value = 9
$ws.Cells.Item(93, 4).Value = "FRANCISCO DIAZ"
$ws.Cells.Item(93, 5).Value = "Walfred Lira"
$ws.Cells.Item(93, 6).Value = "ELECTRICIDAD"
$ws.Cells.Item(93, 7).Value = "Bosquemar"
$ws.Cells.Item(93, 8).Value = "'15"
$ws.Cells.Item(93, 9).Value = "'3"
$ws.Cells.Item(93, 10).Value = "Iluminación"
$ws.Cells.Item(93, 11).Value = "Paused"
$ws.Cells.Item(93, 12).Value = 4
$ws.Cells.Item(93, 13).Value = "L1"
$ws.Cells.Item(93, 14).Value = "2024-07-13 15:56"
$ws.Cells.Item(93, 15).Value = "Fin del día automático"

# Row 94
$ws.Cells.Item(94, 1).Value = "76292e9b-9f21-40d4-bbf0-669d245975c2"
$ws.Cells.Item(94, 2).Value = "2024-07-12 07:53"
$ws.Cells.Item(94, 3).Value = 7
$ws.Cells.Item(94, 4).Value = "CESAR VILLARROEL"
$ws.Cells.Item(94, 5).Value = "Luis Recabal"
$ws.Cells.Item(94, 6).Value = "CARPINTERIA"
$ws.Cells.Item(94, 7).Value = "Puyaral"
$ws.Cells.Item(94, 8).Value = "'15"
$ws.Cells.Item(94, 9).Value = "'2"
$ws.Cells.Item(94, 10).Value = "Barrera de humedad Volcanwrap exterior muros"
$ws.Cells.Item(94, 11).Value = "Finished"
$ws.Cells.Item(94, 12).Value = 4
$ws.Cells.Item(94, 13).Value = "L1"
$ws.Cells.Item(94, 20).Value = "2024-07-12 07:53"

# Row 95
$ws.Cells.Item(95, 1).Value = "d887c85e-f175-4bd4-aa3b-76cf773a7244"
$ws.Cells.Item(95, 2).Value = "2024-07-12 08:57"
$ws.Cells.Item(95, 3).Value = 7
$ws.Cells.Item(95, 4).Value = "CESAR VILLARROEL"
$ws.Cells.Item(95, 5).Value = "Luis Recabal"
$ws.Cells.Item(95, 6).Value = "CARPINTERIA"
$ws.Cells.Item(95, 7).Value = "Puyaral"
$ws.Cells.Item(95, 8).Value = "'15"
$ws.Cells.Item(95, 9).Value = "'2"
$ws.Cells.Item(95, 10).Value = "Piso SPC"
$ws.Cells.Item(95, 11).Value = "Finished"
$ws.Cells.Item(95, 12).Value = 4
$ws.Cells.Item(95, 13).Value = "L1"
$ws.Cells.Item(95, 20).Value = "2024-07-12 08:57"

# Row 96
$ws.Cells.Item(96, 1).Value = "f151ce73-54b7-46ae-b271-7a5dd9b9aefe"
$ws.Cells.Item(96, 2).Value = "2024-07-12 12:39"
$ws.Cells.Item(96, 3).Value = 35
$ws.Cells.Item(96, 4).Value = "FRANCISCO DIAZ"
$ws.Cells.Item(96, 5).Value = "CLAUDIO ROJAS"
$ws.Cells.Item(96, 6).Value = "PINTURA"
$ws.Cells.Item(96, 7).Value = "Puyaral"
$ws.Cells.Item(96, 8).Value = "'14"
$ws.Cells.Item(96, 9).Value = "'2"
$ws.Cells.Item(96, 10).Value = "Pintura Interior (2° mano)"
$ws.Cells.Item(96, 11).Value = "Paused"
$ws.Cells.Item(96, 12).Value = 1
$ws.Cells.Item(96, 13).Value = "L1"
$ws.Cells.Item(96, 14).Value = "2024-07-12 12:39"
$ws.Cells.Item(96, 15).Value = "Fin del día"

# Row 97
$ws.Cells.Item(97, 1).Value = "efcab23b-415e-4d2a-aadf-1f62c0d44829"
$ws.Cells.Item(97, 2).Value = "2024-07-12 13:03"
$ws.Cells.Item(97, 3).Value = 25
$ws.Cells.Item(97, 4).Value = "CESAR VILLARROEL"
$ws.Cells.Item(97, 5).Value = "Celso Martinez"
$ws.Cells.Item(97, 6).Value = "CARPINTERIA"
$ws.Cells.Item(97, 7).Value = "Las Bandurrias"
$ws.Cells.Item(97, 8).Value = "'50"
$ws.Cells.Item(97, 9).Value = "'2"
$ws.Cells.Item(97, 10).Value = "Instalación OSB"
$ws.Cells.Item(97, 11).Value = "Finished"
$ws.Cells.Item(97, 12).Value = 4
$ws.Cells.Item(97, 13).Value = "L2"
$ws.Cells.Item(97, 20).Value = "2024-07-12 13:03"

# Row 98
$ws.Cells.Item(98, 1).Value = "73fb3577-6673-4829-a8a6-5fb9e019701a"
$ws.Cells.Item(98, 2).Value = "2024-07-12 13:15"
$ws.Cells.Item(98, 3).Value = 25
$ws.Cells.Item(98, 4).Value = "CESAR VILLARROEL"
$ws.Cells.Item(98, 5).Value = "Celso Martinez"
$ws.Cells.Item(98, 6).Value = "CARPINTERIA"
$ws.Cells.Item(98, 7).Value = "Puyaral"
$ws.Cells.Item(98, 8).Value = "'50"
$ws.Cells.Item(98, 9).Value = "'2"
$ws.Cells.Item(98, 10).Value = "Guardapolvos y pilastras"
$ws.Cells.Item(98, 11).Value = "Paused"
$ws.Cells.Item(98, 12).Value = 4
$ws.Cells.Item(98, 13).Value = "L2"
$ws.Cells.Item(98, 14).Value = "2024-07-12 16:28"
$ws.Cells.Item(98, 15).Value = "Fin del día"

# Row 99
$ws.Cells.Item(99, 1).Value = "95888789-eb8c-483e-9456-6a5677581cba"
$ws.Cells.Item(99, 2).Value = "2024-07-12 16:15"
$ws.Cells.Item(99, 3).Value = 39
$ws.Cells.Item(99, 4).Value = "FRANCISCO DIAZ"
$ws.Cells.Item(99, 5).Value = "DIEGO RIOSECO"
$ws.Cells.Item(99, 6).Value = "GASFITERÍA"
$ws.Cells.Item(99, 7).Value = "Las Bandurrias"
$ws.Cells.Item(99, 8).Value = "'25"
$ws.Cells.Item(99, 9).Value = "'2"
$ws.Cells.Item(99, 10).Value = "Artefactos sanitarios"
$ws.Cells.Item(99, 11).Value = "Paused"
$ws.Cells.Item(99, 12).Value = 3
$ws.Cells.Item(99, 13).Value = "L2"
$ws.Cells.Item(99, 14).Value = "2024-07-12 16:16"
$ws.Cells.Item(99, 15).Value = "Fin del día"

# Row 100
$ws.Cells.Item(100, 1).Value = "9c10e9aa-0664-48e2-af1f-77c612e2646e"
$ws.Cells.Item(100, 2).Value = "2024-07-12 16:46"
$ws.Cells.Item(100, 3).Value = 25
$ws.Cells.Item(100, 4).Value = "CESAR VILLARROEL"
$ws.Cells.Item(100, 5).Value = "CELSO MARTINEZ"
$ws.Cells.Item(100, 6).Value = "CARPINTERIA"
$ws.Cells.Item(100, 7).Value = "Las Bandurrias"
$ws.Cells.Item(100, 8).Value = "'25"
$ws.Cells.Item(100, 9).Value = "'2"
$ws.Cells.Item(100, 10).Value = "Instalación Escalera"
$ws.Cells.Item(100, 11).Value = "Finished"
$ws.Cells.Item(100, 12).Value = 3
$ws.Cells.Item(100, 13).Value = "L2"
$ws.Cells.Item(100, 20).Value = "2024-07-12 16:47"

# Row 101
$ws.Cells.Item(101, 1).Value = "2433e744-e529-44df-914e-921d5cd4fa25"
$ws.Cells.Item(101, 2).Value = "2024-07-12 16:53"
$ws.Cells.Item(101, 3).Value = 46
$ws.Cells.Item(101, 4).Value = "CAMILO CASTILLO"
$ws.Cells.Item(101, 5).Value = "ABRAHAM BECERRA"
$ws.Cells.Item(101, 6).Value = "PINTURA"
$ws.Cells.Item(101, 7).Value = "Las Bandurrias"
$ws.Cells.Item(101, 8).Value = "'25"
$ws.Cells.Item(101, 9).Value = "'2"
$ws.Cells.Item(101, 10).Value = "Pintura Interior (2° mano)"
$ws.Cells.Item(101, 11).Value = "Finished"
$ws.Cells.Item(101, 12).Value = 3
$ws.Cells.Item(101, 13).Value = "L2"
$ws.Cells.Item(101, 20).Value = "2024-07-12 16:53"
